# TB-383 export fixture regeneration:
#   00162 - 1080p bandwidth  ->  00085 - digital interface
# Renames both worksheets, updates the title cell (A1) on each sheet, and
# reapplies the column widths that Excel recalculated when it rebuilt the
# export (auto-fit drift from the new, differently-sized title text).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "Project_budget_00085_digital in"
$ws2.Name = "Partner_budget_00085_digital in"

$newTitle = "00085 - digital interface - V1.0 - 2023/06/20 - 11:27:54"
$ws1.Range("A1").Value = $newTitle
$ws2.Range("A1").Value = $newTitle

# Column widths (as they end up stored in the sheet XML) for the
# "Project budget" sheet (29 columns) and the "Partner budget" sheet
# (25 columns). Values passed to ColumnWidth are pre-compensated so the
# resulting stored width matches the exported template.
$sheet1Widths = @(50.833333333333336, 6.0, 24.166666666666668, 42.333333333333336, 33.0, 7.5, 14.833333333333334, 17.666666666666668, 14.166666666666666, 14.333333333333334, 34.833333333333336, 38.0, 31.333333333333332, 33.166666666666664, 35.0, 35.0, 35.166666666666664, 38.833333333333336, 38.833333333333336, 18.666666666666668, 23.833333333333332, 20.0, 28.833333333333332, 32.5, 32.5, 10.833333333333334, 9.833333333333334, 45.0, 19.166666666666668)

$sheet2Widths = @(50.833333333333336, 6.0, 24.166666666666668, 42.333333333333336, 33.0, 7.5, 7.0, 7.0, 25.333333333333332, 42.0, 41.5, 17.0, 10.666666666666666, 20.666666666666668, 10.833333333333334, 9.5, 40.166666666666664, 10.666666666666666, 12.5, 10.833333333333334, 7.666666666666667, 7.666666666666667, 7.666666666666667, 7.166666666666667, 8.333333333333334)

for ($i = 0; $i -lt $sheet1Widths.Length; $i++) {
    $ws1.Columns.Item($i + 1).ColumnWidth = $sheet1Widths[$i]
}

for ($i = 0; $i -lt $sheet2Widths.Length; $i++) {
    $ws2.Columns.Item($i + 1).ColumnWidth = $sheet2Widths[$i]
}
